# Generate Report for Handoff
# Refresh the "Latest Handoff Date/Datetime" timestamp for every row that is
# still sitting in "Ready for handoff" (rows 10-16) as well as the row whose
# handback transform failed (row 7), on all three sheets, to reflect the
# newly (re-)generated handoff report timestamp.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

foreach ($r in $rows) {
    $overview.Range("D$r").Value = "2016-21-12 14:21:52"
    $zhcn.Range("E$r").Value = "2016-03-12 14:21:48"
    $dede.Range("E$r").Value = "2016-03-12 14:21:52"
}
